# Using R2Jags instead of rjags for Gibbs-Sampling
#
# DATA_SET sheet: correct the sampling timestamps for patient PAT0001
# (rows 2-3) and append the remaining dosing/sampling events (rows 4-8)
# that belong to the same 2019-12-15/16 admission.
#
# PATIENT / PATHOGEN sheets: update the patient's weight/creatinine
# clearance and the pathogen MIC used by the model.

$wb = $excel.ActiveWorkbook

# Helper: write a value as plain text, bypassing Excel's "looks like a
# date/number" auto-conversion, and leave the cell's style untouched
# (no explicit numeric/text format lingers on the cell afterwards).
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$ws1 = $wb.Worksheets.Item("DATA_SET")

# Row data for the DATA_SET table (rows 2-3 are edits of existing rows,
# rows 4-8 are newly appended rows). Every column except EVID (E) is
# stored as text in this sheet, even when the text looks numeric.
$rows = @(
    @{ Row = 2; A = "2019-12-15"; B = "07:23:00"; C = "1000"; D = "."  ; E = 1; F = "60"; G = "."; H = "."; I = "." },
    @{ Row = 3; A = "2019-12-15"; B = "12:10:00"; C = "."   ; D = "12" ; E = 0; F = "." ; G = "."; H = "."; I = "." },
    @{ Row = 4; A = "2019-12-15"; B = "18:10:00"; C = "."   ; D = "9"  ; E = 0; F = "." ; G = "."; H = "."; I = "." },
    @{ Row = 5; A = "2019-12-15"; B = "19:09:00"; C = "1000"; D = "."  ; E = 1; F = "60"; G = "."; H = "."; I = "." },
    @{ Row = 6; A = "2019-12-16"; B = "06:25:00"; C = "."   ; D = "15" ; E = 0; F = "." ; G = "."; H = "."; I = "." },
    @{ Row = 7; A = "2019-12-16"; B = "07:10:00"; C = "1000"; D = "."  ; E = 1; F = "60"; G = "."; H = "."; I = "." },
    @{ Row = 8; A = "2019-12-16"; B = "12:35:00"; C = "."   ; D = "25" ; E = 0; F = "." ; G = "."; H = "."; I = "." }
)

foreach ($r in $rows) {
    $n = $r.Row
    Set-TextValue $ws1 "A$n" $r.A
    Set-TextValue $ws1 "B$n" $r.B
    Set-TextValue $ws1 "C$n" $r.C
    Set-TextValue $ws1 "D$n" $r.D
    $ws1.Range("E$n").Value = $r.E
    Set-TextValue $ws1 "F$n" $r.F
    Set-TextValue $ws1 "G$n" $r.G
    Set-TextValue $ws1 "H$n" $r.H
    Set-TextValue $ws1 "I$n" $r.I
}

# PATIENT sheet: weight (WT) 70 -> 65, creatinine clearance (CRCL) 120 -> 80
$ws2 = $wb.Worksheets.Item("PATIENT")
$ws2.Range("B2").Value = 65
$ws2.Range("C2").Value = 80

# PATHOGEN sheet: MIC 0.85 -> 5
$ws3 = $wb.Worksheets.Item("PATHOGEN")
$ws3.Range("B2").Value = 5
